$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 236.8597
    3  = 236.5671
    4  = 236.2937
    5  = 237.2014
    6  = 238.6976
    7  = 240.5879
    8  = 241.1502
    9  = 241.1166
    10 = 241.14
    11 = 241.2294
    12 = 241.7053
    13 = 241.2979
    14 = 241.2979
    15 = 241.0549
    16 = 241.1006
    17 = 241.5286
    18 = 241.5213
    19 = 241.3417
    20 = 242.2374
    21 = 241.7493
    22 = 241.49
    23 = 242.0937
    24 = 241.7141
    25 = 241.5181
    26 = 241.9896
    27 = 243.2627
    28 = 241.4055
    29 = 241.0823
    30 = 241.1169
    31 = 240.527
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}
